$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Original layout (before edit) ---
#   Row1: B1..F1 headers
#   Row2: A2=0, B2=40810, C2=<blank>, D2=100, E2=753, F2=753
#   Row3: A3=1, B3=41809, C3=688,    D3=687, E3=756, F3=756
#
# This is rebuilt into a "compare" style report: three grouped columns
# (CandidateNumber / RequisitionNumber / ReferrerPersonNumber), each split
# into HDL vs Data_Val sub-columns, with the row-index column A, and the
# two data records shifted down two rows (with a new third record added).

# --- Propagate the index-column style (bold/boxed/centered) to the new data
#     rows while row 3 still carries it, then wipe row 3 completely ---
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)
$ws.Rows("3").Clear()

# --- Strip the existing per-cell borders from the header block before
#     merging, so Excel does not split the box border across the merged
#     cells (keeps every header cell on the same, single shared style) ---
$ws.Range("A1:G2").ClearFormats()

# --- Merge the group header cells ---
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()

# --- Row 1: merged group headers (CandidateNumber / RequisitionNumber / ReferrerPersonNumber) ---
$ws.Range("A1").Value = ""
$ws.Range("B1").Value = "CandidateNumber"
$ws.Range("D1").Value = "RequisitionNumber"
$ws.Range("F1").Value = "ReferrerPersonNumber"

# --- Row 2: HDL / Data_Val sub-headers ---
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "HDL"
$ws.Range("C2").Value = "Data_Val"
$ws.Range("D2").Value = "HDL"
$ws.Range("E2").Value = "Data_Val"
$ws.Range("F2").Value = "HDL"
$ws.Range("G2").Value = "Data_Val"

# --- Re-apply the original bold/boxed/centered header style uniformly
#     across the whole new header block (same style index for every cell).
#     A4 already carries that style (copied from the old A3 above). ---
$ws.Range("A4").Copy()
$ws.Range("A1:G2").PasteSpecial(-4122)

# --- Row 4 (record 0) ---
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 40810
$ws.Range("C4").Value = 40810
$ws.Range("D4").Value = 753
$ws.Range("E4").Value = 753
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = ""

# --- Row 5 (record 1) ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 41809
$ws.Range("C5").Value = 41809
$ws.Range("D5").Value = 756
$ws.Range("E5").Value = 756
$ws.Range("F5").Value = 688
$ws.Range("G5").Value = 687

# --- Row 6 (record 2, newly added) ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 41809
$ws.Range("C6").Value = 41809
$ws.Range("D6").Value = 756
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 688
$ws.Range("G6").Value = 688
